# Edit script: "Descripción para la plataforma de TFGs actualizada"
$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Title: "APLICACIÓN PARA LA EXPOSICIÓN Y VALORACIÓN DE VIDEOJUEGOS"
#    -> "APLICACIÓN PARA EL INTERCAMBIO DE CONTENIDO ENTRE AFICIONADOS A LOS VIDEOJUEGOS"
# ---------------------------------------------------------------
$d.Content.Find.Execute("APLICACIÓN PARA ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "APLICACIÓN PARA", 2) | Out-Null
$d.Content.Find.Execute("LA EXPOSICIÓN Y VALORACIÓN DE VIDEOJUEGOS", $true, $false, $false, $false, $false, `
    $true, 1, $false, " EL INTERCAMBIO DE CONTENIDO ENTRE AFICIONADOS A LOS VIDEOJUEGOS", 2) | Out-Null

# ---------------------------------------------------------------
# 2) DESCRIPCIÓN paragraph: reword body text
# ---------------------------------------------------------------
$oldDescBody = " Se pretende diseñar e implementar una aplicación basada en el intercambio de información entre jugadores de videojuegos, de forma que puedan describir y explicar diferentes juegos. Además de aportar información sobre nuevos videojuegos, los usuarios podrán valorar los ya propuestos por otros usuarios y añadir comentarios o aportaciones."
$newDescBody = " Se pretende diseñar e implementar una aplicación, que facilite la comunicación y el intercambio de información entre aficionados al mundo de los videojuegos, permitiendo la comunicación directa entre ellos, la publicación y visualización de contenido, la valoración y opinión de dicho contenido, y demás funcionalidades que vayan surgiendo durante el desarrollo del producto."
$d.Content.Find.Execute($oldDescBody, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newDescBody, 2) | Out-Null

# Move/re-create the "_GoBack" bookmark right after "...una aplicación" (before the comma)
$gb = $d.Content
$gb.Find.Execute("Se pretende diseñar e implementar una aplicación", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$gb.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $gb) | Out-Null

# ---------------------------------------------------------------
# 3) "En primer lugar..." paragraph: reword + split into two paragraphs
# ---------------------------------------------------------------
$oldPara2 = "En primer lugar, se realizará una fase de análisis de requisitos y diseño de la aplicación, a través de encuestas a potenciales usuarios y de investigación sobre diferentes frameworks. Después se realizará una implementación basada en prototipos, de modo que se irán produciendo versiones más sencillas al principio, que se irán ampliando con versiones más completas a lo largo del trabajo (según los requisitos iniciales y los que vayan surgiendo de las versiones)."
$newPara2 = "Será necesario realizar una fase inicial de estudio y análisis de los “frameworks” disponibles, así como un análisis competitivo y unas encuestas a futuros usuarios para obtener los requisitos iniciales del sistema."
$d.Content.Find.Execute($oldPara2, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newPara2, 2) | Out-Null

# Insert the new paragraph right after it
$frameworksPara = $d.Content
$frameworksPara.Find.Execute($newPara2, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$frameworksPara.Collapse(0) | Out-Null
$frameworksPara.InsertParagraphAfter() | Out-Null
$frameworksPara.Collapse(0) | Out-Null
$newLifecyclePara = $frameworksPara.Next(4) # wdParagraph = 4, move to the new paragraph
$newLifecyclePara.Text = "Se seguirá un ciclo de vida incremental iterativo, pudiendo obtener nuevos requisitos y cambios, para que se vayan generando prototipos cada vez más ricos en funcionalidad."
